# Horarios actualizados Linea 141 - 819
# Applies the scraped-data refresh across the three worksheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: LP1912
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("LP1912")

$ws.Range("A2").Value = "Última actualización: 04:56:11"
$ws.Range("A3").Value = "Total filas: 38"

$ws.Cells.Item(20,1).Value = "04:56:11"
$ws.Cells.Item(20,2).Value = "05:16"
$ws.Cells.Item(20,3).Value = "17_ROMERO"
$ws.Cells.Item(20,4).Value = 20
$ws.Cells.Item(20,5).Value = "LP1912"

$ws.Cells.Item(22,1).Value = "04:56:11"
$ws.Cells.Item(22,2).Value = "05:22"
$ws.Cells.Item(22,3).Value = "23_HERNANDEZ"
$ws.Cells.Item(22,4).Value = 26
$ws.Cells.Item(22,5).Value = "LP1912"

$ws.Cells.Item(25,1).Value = "04:56:11"
$ws.Cells.Item(25,2).Value = "05:35"
$ws.Cells.Item(25,3).Value = "215B_EL PATO"
$ws.Cells.Item(25,4).Value = 39
$ws.Cells.Item(25,5).Value = "LP1912"

$ws.Cells.Item(27,1).Value = "04:56:11"
$ws.Cells.Item(27,2).Value = "05:46"
$ws.Cells.Item(27,3).Value = "15_ABASTO"
$ws.Cells.Item(27,4).Value = 50
$ws.Cells.Item(27,5).Value = "LP1912"

$ws.Cells.Item(28,1).Value = "04:56:11"
$ws.Cells.Item(28,2).Value = "05:54"
$ws.Cells.Item(28,3).Value = "10_OLMOS"
$ws.Cells.Item(28,4).Value = 58
$ws.Cells.Item(28,5).Value = "LP1912"

$ws.Cells.Item(29,1).Value = "04:56:11"
$ws.Cells.Item(29,2).Value = "06:04"
$ws.Cells.Item(29,3).Value = "16_SANTA ANA"
$ws.Cells.Item(29,4).Value = 68
$ws.Cells.Item(29,5).Value = "LP1912"

$ws.Cells.Item(30,1).Value = "04:56:11"
$ws.Cells.Item(30,2).Value = "06:11"
$ws.Cells.Item(30,3).Value = "215A_EL PATO"
$ws.Cells.Item(30,4).Value = 75
$ws.Cells.Item(30,5).Value = "LP1912"

$ws.Cells.Item(31,1).Value = "04:48:57"
$ws.Cells.Item(31,2).Value = "06:13"
$ws.Cells.Item(31,3).Value = "225_HARAS DEL SUR"
$ws.Cells.Item(31,4).Value = 85
$ws.Cells.Item(31,5).Value = "LP1912"

$ws.Cells.Item(32,1).Value = "04:56:11"
$ws.Cells.Item(32,2).Value = "06:14"
$ws.Cells.Item(32,3).Value = "225_HARAS DEL SUR"
$ws.Cells.Item(32,4).Value = 78
$ws.Cells.Item(32,5).Value = "LP1912"

$ws.Cells.Item(33,1).Value = "04:48:57"
$ws.Cells.Item(33,2).Value = "06:20"
$ws.Cells.Item(33,3).Value = "26_HERNANDEZ"
$ws.Cells.Item(33,4).Value = 92
$ws.Cells.Item(33,5).Value = "LP1912"

$ws.Cells.Item(34,1).Value = "04:56:11"
$ws.Cells.Item(34,2).Value = "06:21"
$ws.Cells.Item(34,3).Value = "26_HERNANDEZ"
$ws.Cells.Item(34,4).Value = 85
$ws.Cells.Item(34,5).Value = "LP1912"

$ws.Cells.Item(35,1).Value = "04:48:57"
$ws.Cells.Item(35,2).Value = "06:26"
$ws.Cells.Item(35,3).Value = "23_HERNANDEZ"
$ws.Cells.Item(35,4).Value = 98
$ws.Cells.Item(35,5).Value = "LP1912"

$ws.Cells.Item(36,1).Value = "04:56:11"
$ws.Cells.Item(36,2).Value = "06:27"
$ws.Cells.Item(36,3).Value = "23_HERNANDEZ"
$ws.Cells.Item(36,4).Value = 91
$ws.Cells.Item(36,5).Value = "LP1912"

$ws.Cells.Item(37,1).Value = "04:48:57"
$ws.Cells.Item(37,2).Value = "06:29"
$ws.Cells.Item(37,3).Value = "86_EST CHICA-ESC AGRARIA"
$ws.Cells.Item(37,4).Value = 101
$ws.Cells.Item(37,5).Value = "LP1912"

$ws.Cells.Item(38,1).Value = "04:56:11"
$ws.Cells.Item(38,2).Value = "06:30"
$ws.Cells.Item(38,3).Value = "86_EST CHICA-ESC AGRARIA"
$ws.Cells.Item(38,4).Value = 94
$ws.Cells.Item(38,5).Value = "LP1912"

$ws.Cells.Item(39,1).Value = "04:56:11"
$ws.Cells.Item(39,2).Value = "06:31"
$ws.Cells.Item(39,3).Value = "16_SANTA ANA"
$ws.Cells.Item(39,4).Value = 95
$ws.Cells.Item(39,5).Value = "LP1912"

$ws.Cells.Item(40,1).Value = "04:48:57"
$ws.Cells.Item(40,2).Value = "06:43"
$ws.Cells.Item(40,3).Value = "225_C ROCA-H SUR"
$ws.Cells.Item(40,4).Value = 115
$ws.Cells.Item(40,5).Value = "LP1912"

$ws.Cells.Item(41,1).Value = "04:56:11"
$ws.Cells.Item(41,2).Value = "06:44"
$ws.Cells.Item(41,3).Value = "225_C ROCA-H SUR"
$ws.Cells.Item(41,4).Value = 108
$ws.Cells.Item(41,5).Value = "LP1912"

$ws.Cells.Item(42,1).Value = "04:48:57"
$ws.Cells.Item(42,2).Value = "06:46"
$ws.Cells.Item(42,3).Value = "215C_EL PATO"
$ws.Cells.Item(42,4).Value = 118
$ws.Cells.Item(42,5).Value = "LP1912"

$ws.Cells.Item(43,1).Value = "04:56:11"
$ws.Cells.Item(43,2).Value = "06:47"
$ws.Cells.Item(43,3).Value = "215C_EL PATO"
$ws.Cells.Item(43,4).Value = 111
$ws.Cells.Item(43,5).Value = "LP1912"

# ---------------------------------------------------------------------------
# Sheet 2: LP1912-215
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("LP1912-215")

$ws2.Range("A2").Value = "Última actualización: 04:56:11"
$ws2.Range("A3").Value = "Total filas: 11"

$ws2.Cells.Item(13,1).Value = "04:56:11"
$ws2.Cells.Item(13,2).Value = "05:35"
$ws2.Cells.Item(13,3).Value = "215B_EL PATO"
$ws2.Cells.Item(13,4).Value = 39
$ws2.Cells.Item(13,5).Value = "LP1912"

$ws2.Cells.Item(14,1).Value = "04:56:11"
$ws2.Cells.Item(14,2).Value = "06:11"
$ws2.Cells.Item(14,3).Value = "215A_EL PATO"
$ws2.Cells.Item(14,4).Value = 75
$ws2.Cells.Item(14,5).Value = "LP1912"

$ws2.Cells.Item(16,1).Value = "04:56:11"
$ws2.Cells.Item(16,2).Value = "06:47"
$ws2.Cells.Item(16,3).Value = "215C_EL PATO"
$ws2.Cells.Item(16,4).Value = 111
$ws2.Cells.Item(16,5).Value = "LP1912"

# ---------------------------------------------------------------------------
# Sheet 3: 6203-6173
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("6203-6173")

$ws3.Range("A2").Value = "Última actualización: 04:56:11"

$ws3.Cells.Item(8,1).Value = "04:56:11"
$ws3.Cells.Item(8,2).Value = "05:44"
$ws3.Cells.Item(8,3).Value = "215A_LA PLATA"
$ws3.Cells.Item(8,4).Value = 48
$ws3.Cells.Item(8,5).Value = "L6173"

$ws3.Cells.Item(10,1).Value = "04:56:11"
$ws3.Cells.Item(10,2).Value = "06:09"
$ws3.Cells.Item(10,3).Value = "215A_LA PLATA"
$ws3.Cells.Item(10,4).Value = 73
$ws3.Cells.Item(10,5).Value = "L6173"

$ws3.Cells.Item(12,1).Value = "04:56:11"
$ws3.Cells.Item(12,2).Value = "06:33"
$ws3.Cells.Item(12,3).Value = "215C_LA PLATA"
$ws3.Cells.Item(12,4).Value = 97
$ws3.Cells.Item(12,5).Value = "L6203"
